$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 - Gross Margin
$ws.Range("D15").Value = 1.0028
$ws.Range("E15").Value = 0.9938
$ws.Range("F15").Value = 0.9959
$ws.Range("G15").Value = 1.0

# Row 16 - EBIT Margin
$ws.Range("D16").Value = 0.1861
$ws.Range("E16").Value = 0.1866
$ws.Range("F16").Value = 0.1949
$ws.Range("G16").Value = 0.2013

# Row 17 - EBT margin
$ws.Range("D17").Value = 0.1292
$ws.Range("E17").Value = 0.1409
$ws.Range("F17").Value = 0.1565
$ws.Range("G17").Value = 0.1746

# Row 18 - Net Profit Margin
$ws.Range("D18").Value = 0.1003
$ws.Range("E18").Value = 0.1068
$ws.Range("F18").Value = 0.1158
$ws.Range("G18").Value = 0.1308

# Row 19 - Free Cash Flow Margin
$ws.Range("D19").Value = 0.4797
$ws.Range("E19").Value = 0.3466
$ws.Range("F19").Value = 0.3807
$ws.Range("G19").Value = 0.0935

# Row 27 - EBITDA Margin
$ws.Range("D27").Value = 0.2077
$ws.Range("E27").Value = 0.2076
$ws.Range("F27").Value = 0.2151
$ws.Range("G27").Value = 0.2181

# Row 28 - Operating Cash Flow Margin
$ws.Range("D28").Value = 0.4949
$ws.Range("E28").Value = 0.3629
$ws.Range("F28").Value = 0.3989
$ws.Range("G28").Value = 0.1119
